$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 123816.164
$ws.Range("I6").Value = 185287.75
$ws.Range("J6").Value = 873
$ws.Range("K6").Value = 555863.25
$ws.Range("L6").Value = 2619
$ws.Range("M6").Value = -555751.25
$ws.Range("N6").Value = -2843

$ws.Range("H17").Value = 345.93332
$ws.Range("J17").Value = 345.93332
$ws.Range("L17").Value = 1037.79996
$ws.Range("N17").Value = -1373.79996

$ws.Range("H28").Value = 1259.8889
$ws.Range("I28").Value = 689.8333
$ws.Range("J28").Value = 2400
$ws.Range("K28").Value = 689.8333
$ws.Range("L28").Value = 2400
$ws.Range("M28").Value = -204.8333
$ws.Range("N28").Value = -3370

$ws.Range("H70").Value = 1605.1
$ws.Range("I70").Value = 1999.75
$ws.Range("K70").Value = 5999.25
$ws.Range("M70").Value = -5729.25

$ws.Range("H73").Value = 1605.1
$ws.Range("I73").Value = 1999.75
$ws.Range("K73").Value = 5999.25
$ws.Range("M73").Value = -5063.25

$ws.Range("H125").Value = 2016
$ws.Range("I125").Value = 2241.5
$ws.Range("J125").Value = 1852
$ws.Range("K125").Value = 20173.5
$ws.Range("L125").Value = 16668
$ws.Range("M125").Value = -17713.5
$ws.Range("N125").Value = -21588

$ws.Range("H132").Value = 5686958.5
$ws.Range("I132").Value = 6103011.5
$ws.Range("J132").Value = 900
$ws.Range("K132").Value = 18309034.5
$ws.Range("L132").Value = 2700
$ws.Range("M132").Value = -18306504.5
$ws.Range("N132").Value = -7760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 1857.7142
$ws.Range("I10").Value = 1004
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 1004
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = -834
$ws.Range("N10").Value = -2340

$ws.Range("H132").Value = 2456.75
$ws.Range("I132").Value = 2100.7693
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 6302.3079
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -3772.3079
$ws.Range("N132").Value = -17057.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H99").Value = 2092.5
$ws.Range("I99").Value = 1908
$ws.Range("K99").Value = 1908
$ws.Range("M99").Value = -410

$ws.Range("H134").Value = 3064
$ws.Range("I134").Value = 3064
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9192
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -6657

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 407
$ws.Range("I10").Value = 407
$ws.Range("K10").Value = 407
$ws.Range("M10").Value = -268

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 333.82352
$ws.Range("I18").Value = 294.8
$ws.Range("J18").Value = 626.5
$ws.Range("K18").Value = 884.4000000000001
$ws.Range("L18").Value = 1879.5
$ws.Range("M18").Value = -715.4000000000001
$ws.Range("N18").Value = -2217.5

$ws.Range("H92").Value = 1200
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()

$ws.Range("H107").Value = 1391138.8
$ws.Range("J107").Value = 1738435.9
$ws.Range("L107").Value = 5215307.699999999
$ws.Range("N107").Value = -5219147.699999999

$ws.Range("H129").Value = 22738542
$ws.Range("I129").Value = 83333650
$ws.Range("J129").Value = 15375
$ws.Range("K129").Value = 250000950
$ws.Range("L129").Value = 46125
$ws.Range("M129").Value = -249995950
$ws.Range("N129").Value = -56125

$ws.Range("H131").Value = 863.7
$ws.Range("I131").Value = 710
$ws.Range("J131").Value = 868.4536000000001
$ws.Range("K131").Value = 2130
$ws.Range("L131").Value = 2605.3608
$ws.Range("M131").Value = 2910
$ws.Range("N131").Value = -12685.3608

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 5858500
$ws.Range("J12").Value = 4833.3335
$ws.Range("L12").Value = 4833.3335
$ws.Range("N12").Value = -5113.3335

$ws.Range("H70").Value = 81436.5
$ws.Range("I70").Value = 147011.14
$ws.Range("J70").Value = 4932.75
$ws.Range("K70").Value = 147011.14
$ws.Range("L70").Value = 4932.75
$ws.Range("M70").Value = -146741.14
$ws.Range("N70").Value = -5472.75

$ws.Range("H73").Value = 81436.5
$ws.Range("I73").Value = 147011.14
$ws.Range("J73").Value = 4932.75
$ws.Range("K73").Value = 147011.14
$ws.Range("L73").Value = 4932.75
$ws.Range("M73").Value = -146075.14
$ws.Range("N73").Value = -6804.75

$ws.Range("H126").Value = 2377.4736
$ws.Range("I126").Value = 2418.4
$ws.Range("J126").Value = 2362.8572
$ws.Range("K126").Value = 7255.200000000001
$ws.Range("L126").Value = 7088.571599999999
$ws.Range("M126").Value = -4785.200000000001
$ws.Range("N126").Value = -12028.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 6000
$ws.Range("J24").Value = 10000
$ws.Range("L24").Value = 10000
$ws.Range("N24").Value = -10686

$ws.Range("H40").Value = 85583.25
$ws.Range("J40").Value = 2499.9
$ws.Range("L40").Value = 2499.9
$ws.Range("N40").Value = -2771.9

$ws.Range("H61").Value = 4296.5
$ws.Range("I61").Value = 3874.75
$ws.Range("J61").Value = 5140
$ws.Range("K61").Value = 3874.75
$ws.Range("L61").Value = 5140
$ws.Range("M61").Value = -3672.75
$ws.Range("N61").Value = -5544

$ws.Range("H113").Value = 4296.5
$ws.Range("I113").Value = 3874.75
$ws.Range("J113").Value = 5140
$ws.Range("K113").Value = 3874.75
$ws.Range("L113").Value = 5140
$ws.Range("M113").Value = -1704.75
$ws.Range("N113").Value = -9480

$ws.Range("H136").Value = 2100.5
$ws.Range("I136").Value = 1999.8334
$ws.Range("J136").Value = 2402.5
$ws.Range("K136").Value = 5999.5002
$ws.Range("L136").Value = 7207.5
$ws.Range("M136").Value = -3449.5002
$ws.Range("N136").Value = -12307.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 162356.12
$ws.Range("I107").Value = 57069.8
$ws.Range("J107").Value = 337833.34
$ws.Range("K107").Value = 171209.4
$ws.Range("L107").Value = 1013500.02
$ws.Range("M107").Value = -169289.4
$ws.Range("N107").Value = -1017340.02

$ws.Range("H132").Value = 3217.9614
$ws.Range("I132").Value = 3387.4443
$ws.Range("J132").Value = 2836.625
$ws.Range("K132").Value = 10162.3329
$ws.Range("L132").Value = 8509.875
$ws.Range("M132").Value = -7632.332900000001
$ws.Range("N132").Value = -13569.875

$ws.Range("H136").Value = 1190.2
$ws.Range("I136").Value = 470.2903
$ws.Range("J136").Value = 2784.2856
$ws.Range("K136").Value = 1410.8709
$ws.Range("L136").Value = 8352.856800000001
$ws.Range("M136").Value = 1139.1291
$ws.Range("N136").Value = -13452.8568

